$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header cell A1: "names" -> "name"
$ws.Range("A1").Value = "name"

# 2. Right-align the other header cells (B1:D1) -- creates a new cellXf
#    -4152 == xlRight
$ws.Range("B1:D1").HorizontalAlignment = -4152

# 3. Freeze panes at B2 (freeze first row + first column)
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 4. Column widths (best-fit-like sizing)
$ws.Columns.Item(1).ColumnWidth = 5.1666666666666667
$ws.Columns.Item(2).ColumnWidth = 7.3333333333333333
$ws.Columns.Item(3).ColumnWidth = 7.3333333333333333
$ws.Columns.Item(4).ColumnWidth = 7.3333333333333333
